$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This string is shared by the Overview summary sheet (columns B/C) as
#    well as the per-language "Status" column (C) on the zh-cn/de-de sheets.
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Latest Handback DateTime (column H) now has real timestamps instead of
#    the "0001-01-01 00:00:00" placeholder - one new timestamp per language.
# ---------------------------------------------------------------------------
$wsZh.Range("H2").Value = "2016-03-14 03:12:49"
$wsZh.Range("H3").Value = "2016-03-14 03:12:49"
$wsDe.Range("H2").Value = "2016-03-14 03:12:54"
$wsDe.Range("H3").Value = "2016-03-14 03:12:54"

# ---------------------------------------------------------------------------
# 3. Populate the new "Latest Target File" (F) and "Latest Handback File" (G)
#    columns for both rows of the zh-cn and de-de sheets. Values mirror the
#    existing "Source File Name" (A) / "Latest Handoff File" (D) cells, and
#    carry the same hyperlink + visual style as those source cells.
# ---------------------------------------------------------------------------
function Set-HandbackColumns {
    param($ws, $row, $sourceDisplay, $sourceUrl, $handoffDisplay, $handoffUrl)

    $fCell = $ws.Range("F" + $row)
    $fCell.Value = $sourceDisplay
    $fCell.Font.Name = "Calibri"
    $fCell.Font.Size = 11
    $fCell.Font.Underline = $true
    $fCell.Font.Color = 15570276
    $ws.Hyperlinks.Add($fCell, $sourceUrl, "", "", $sourceDisplay) | Out-Null

    $gCell = $ws.Range("G" + $row)
    $gCell.Value = $handoffDisplay
    $gCell.Font.Name = "Calibri"
    $gCell.Font.Size = 11
    $gCell.Font.Underline = $true
    $gCell.Font.Color = 15570276
    $ws.Hyperlinks.Add($gCell, $handoffUrl, "", "", $handoffDisplay) | Out-Null
}

$mdUrlBase = "https://github.com/OpenLocalizationTest/oltest/blob/2fbd0cd0a0924a31bd384d0212558b9e1c44234c/e2e/"

# zh-cn sheet, row 2 (0ef960b6-...)
Set-HandbackColumns $wsZh 2 `
    "0ef960b6-5170-4df7-a2be-719d30bb3004.md" `
    ($mdUrlBase + "0ef960b6-5170-4df7-a2be-719d30bb3004.md") `
    "0ef960b6-5170-4df7-a2be-719d30bb3004.e80a8d7ea7f92892458ae97b458a5ccba6ec39b1.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/36e0201b89555e0ed0f467310214338f356fc762/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/0ef960b6-5170-4df7-a2be-719d30bb3004.e80a8d7ea7f92892458ae97b458a5ccba6ec39b1.zh-cn.xlf"

# zh-cn sheet, row 3 (d54eafbd-...)
Set-HandbackColumns $wsZh 3 `
    "d54eafbd-4be1-46c6-8533-b07a3286e8cf.md" `
    ($mdUrlBase + "d54eafbd-4be1-46c6-8533-b07a3286e8cf.md") `
    "d54eafbd-4be1-46c6-8533-b07a3286e8cf.66f9f073c0b619251c8150db4e1f2be97f419b26.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/36e0201b89555e0ed0f467310214338f356fc762/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/d54eafbd-4be1-46c6-8533-b07a3286e8cf.66f9f073c0b619251c8150db4e1f2be97f419b26.zh-cn.xlf"

# de-de sheet, row 2 (0ef960b6-...)
Set-HandbackColumns $wsDe 2 `
    "0ef960b6-5170-4df7-a2be-719d30bb3004.md" `
    ($mdUrlBase + "0ef960b6-5170-4df7-a2be-719d30bb3004.md") `
    "0ef960b6-5170-4df7-a2be-719d30bb3004.e80a8d7ea7f92892458ae97b458a5ccba6ec39b1.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e4581585cf8373c8acfba6f3e869a0b1c8a52727/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/0ef960b6-5170-4df7-a2be-719d30bb3004.e80a8d7ea7f92892458ae97b458a5ccba6ec39b1.de-de.xlf"

# de-de sheet, row 3 (d54eafbd-...)
Set-HandbackColumns $wsDe 3 `
    "d54eafbd-4be1-46c6-8533-b07a3286e8cf.md" `
    ($mdUrlBase + "d54eafbd-4be1-46c6-8533-b07a3286e8cf.md") `
    "d54eafbd-4be1-46c6-8533-b07a3286e8cf.66f9f073c0b619251c8150db4e1f2be97f419b26.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e4581585cf8373c8acfba6f3e869a0b1c8a52727/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/d54eafbd-4be1-46c6-8533-b07a3286e8cf.66f9f073c0b619251c8150db4e1f2be97f419b26.de-de.xlf"

Write-Host "Handback report generated."
